$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 13.57210581069602
$ws.Cells.Item(2, 3).Value = 6.515323073507815
$ws.Cells.Item(2, 4).Value = 4.03651090965529
$ws.Cells.Item(2, 5).Value = 7.096331533348716
$ws.Cells.Item(2, 6).Value = 30.36538958399375
$ws.Cells.Item(2, 7).Value = 39.25442742555294
$ws.Cells.Item(2, 9).Value = 4.6767542210036
$ws.Cells.Item(2, 10).Value = 13.2050625099819
$ws.Cells.Item(2, 11).Value = 21.60859420884773
$ws.Cells.Item(2, 12).Value = 5.970802686893333
$ws.Cells.Item(2, 13).Value = 11.10766169490169
$ws.Cells.Item(2, 14).Value = 6.927698423236116
$ws.Cells.Item(3, 2).Value = 12.75007401937467
$ws.Cells.Item(3, 3).Value = 6.15474171981099
$ws.Cells.Item(3, 4).Value = 3.9037908453589
$ws.Cells.Item(3, 5).Value = 7.046483321537263
$ws.Cells.Item(3, 6).Value = 30.03641002603807
$ws.Cells.Item(3, 7).Value = 38.73790592108613
$ws.Cells.Item(3, 9).Value = 4.901066898947545
$ws.Cells.Item(3, 10).Value = 13.194567040139
$ws.Cells.Item(3, 11).Value = 21.57278157321411
$ws.Cells.Item(3, 12).Value = 5.993547212795233
$ws.Cells.Item(3, 13).Value = 10.39832711352476
$ws.Cells.Item(3, 14).Value = 6.820866144484122
$ws.Cells.Item(4, 2).Value = 12.21949320418245
$ws.Cells.Item(4, 3).Value = 5.925574704564259
$ws.Cells.Item(4, 4).Value = 3.822369197966519
$ws.Cells.Item(4, 5).Value = 7.017619996254957
$ws.Cells.Item(4, 6).Value = 29.84300988690177
$ws.Cells.Item(4, 7).Value = 38.43086898558787
$ws.Cells.Item(4, 9).Value = 5.044121308339062
$ws.Cells.Item(4, 10).Value = 13.19173416981502
$ws.Cells.Item(4, 11).Value = 21.55536881246812
$ws.Cells.Item(4, 12).Value = 6.008138471368881
$ws.Cells.Item(4, 13).Value = 9.938484581097121
$ws.Cells.Item(4, 14).Value = 6.758496800405826
$ws.Cells.Item(5, 2).Value = 11.99524872621168
$ws.Cells.Item(5, 3).Value = 5.836289181895888
$ws.Cells.Item(5, 4).Value = 3.791661641544077
$ws.Cells.Item(5, 5).Value = 7.005832783243787
$ws.Cells.Item(5, 6).Value = 29.75292435547073
$ws.Cells.Item(5, 7).Value = 38.28493288044617
$ws.Cells.Item(5, 9).Value = 5.106116828719088
$ws.Cells.Item(5, 10).Value = 13.18707843047154
$ws.Cells.Item(5, 11).Value = 21.5395419655166
$ws.Cells.Item(5, 12).Value = 6.013629370414431
$ws.Cells.Item(5, 13).Value = 9.746388067038351
$ws.Cells.Item(5, 14).Value = 6.735560174685337
$ws.Cells.Item(6, 2).Value = 11.95562370515034
$ws.Cells.Item(6, 3).Value = 5.828751227398081
$ws.Cells.Item(6, 4).Value = 3.789510895340663
$ws.Cells.Item(6, 5).Value = 7.003330565858317
$ws.Cells.Item(6, 6).Value = 29.72166230874542
$ws.Cells.Item(6, 7).Value = 38.23226765274237
$ws.Cells.Item(6, 9).Value = 5.119461669241001
$ws.Cells.Item(6, 10).Value = 13.18099551157908
$ws.Cells.Item(6, 11).Value = 21.52494109723868
$ws.Cells.Item(6, 12).Value = 6.013800617030649
$ws.Cells.Item(6, 13).Value = 9.715850131863755
$ws.Cells.Item(6, 14).Value = 6.733823227692494
$ws.Cells.Item(7, 2).Value = 12.21110591187978
$ws.Cells.Item(7, 3).Value = 5.944214317608671
$ws.Cells.Item(7, 4).Value = 3.829886857042802
$ws.Cells.Item(7, 5).Value = 7.015899253351378
$ws.Cells.Item(7, 6).Value = 29.79710042703192
$ws.Cells.Item(7, 7).Value = 38.35119260267115
$ws.Cells.Item(7, 9).Value = 5.052827668038702
$ws.Cells.Item(7, 10).Value = 13.17705701592333
$ws.Cells.Item(7, 11).Value = 21.5223358260594
$ws.Cells.Item(7, 12).Value = 6.006172936278931
$ws.Cells.Item(7, 13).Value = 9.940532507648939
$ws.Cells.Item(7, 14).Value = 6.763665738796508
$ws.Cells.Item(8, 2).Value = 13.28736734349009
$ws.Cells.Item(8, 3).Value = 6.417237562291716
$ws.Cells.Item(8, 4).Value = 4.00079304625284
$ws.Cells.Item(8, 5).Value = 7.076727747938505
$ws.Cells.Item(8, 6).Value = 30.19206008542844
$ws.Cells.Item(8, 7).Value = 38.97340580930314
$ws.Cells.Item(8, 9).Value = 4.763002327891014
$ws.Cells.Item(8, 10).Value = 13.1815257147065
$ws.Cells.Item(8, 11).Value = 21.55223432838217
$ws.Cells.Item(8, 12).Value = 5.97585983867434
$ws.Cells.Item(8, 13).Value = 10.87376777450675
$ws.Cells.Item(8, 14).Value = 6.897280420869863
$ws.Cells.Item(9, 2).Value = 15.19816956125358
$ws.Cells.Item(9, 3).Value = 7.250636929345772
$ws.Cells.Item(9, 4).Value = 4.321111823489655
$ws.Cells.Item(9, 5).Value = 7.212458480519805
$ws.Cells.Item(9, 6).Value = 31.11264405956542
$ws.Cells.Item(9, 7).Value = 40.40824936606447
$ws.Cells.Item(9, 9).Value = 4.222991148332021
$ws.Cells.Item(9, 10).Value = 13.24583124206785
$ws.Cells.Item(9, 11).Value = 21.70663855509147
$ws.Cells.Item(9, 12).Value = 5.924896407712259
$ws.Cells.Item(9, 13).Value = 12.50429692042874
$ws.Cells.Item(9, 14).Value = 7.174580935450442
$ws.Cells.Item(10, 2).Value = 16.46954784874517
$ws.Cells.Item(10, 3).Value = 7.83368825221585
$ws.Cells.Item(10, 4).Value = 4.556495316305555
$ws.Cells.Item(10, 5).Value = 7.319501809027967
$ws.Cells.Item(10, 6).Value = 31.80897001449433
$ws.Cells.Item(10, 7).Value = 41.47304914063869
$ws.Cells.Item(10, 9).Value = 3.854583232053166
$ws.Cells.Item(10, 10).Value = 13.3056641635775
$ws.Cells.Item(10, 11).Value = 21.82928022995617
$ws.Cells.Item(10, 12).Value = 5.889401687492829
$ws.Cells.Item(10, 13).Value = 13.58829540736649
$ws.Cells.Item(10, 14).Value = 7.394761942908413
$ws.Cells.Item(11, 2).Value = 17.01253010930248
$ws.Cells.Item(11, 3).Value = 8.115497166571481
$ws.Cells.Item(11, 4).Value = 4.67286204341846
$ws.Cells.Item(11, 5).Value = 7.367413322334845
$ws.Cells.Item(11, 6).Value = 32.06806930452271
$ws.Cells.Item(11, 7).Value = 41.85328255187489
$ws.Cells.Item(11, 9).Value = 3.70604705462508
$ws.Cells.Item(11, 10).Value = 13.31467686221696
$ws.Cells.Item(11, 11).Value = 21.83974689238569
$ws.Cells.Item(11, 12).Value = 5.870716147699884
$ws.Cells.Item(11, 13).Value = 14.06206470338538
$ws.Cells.Item(11, 14).Value = 7.505587495547651
$ws.Cells.Item(12, 2).Value = 17.22198475356279
$ws.Cells.Item(12, 3).Value = 8.207514177093675
$ws.Cells.Item(12, 4).Value = 4.710989485543037
$ws.Cells.Item(12, 5).Value = 7.387121220899605
$ws.Cells.Item(12, 6).Value = 32.20329318385911
$ws.Cells.Item(12, 7).Value = 42.06033671611956
$ws.Cells.Item(12, 9).Value = 3.643359817439024
$ws.Cells.Item(12, 10).Value = 13.3309873860967
$ws.Cells.Item(12, 11).Value = 21.87201433520336
$ws.Cells.Item(12, 12).Value = 5.86541409751525
$ws.Cells.Item(12, 13).Value = 14.23494376807113
$ws.Cells.Item(12, 14).Value = 7.543652834114553
$ws.Cells.Item(13, 2).Value = 17.17615293963311
$ws.Cells.Item(13, 3).Value = 8.185088996629716
$ws.Cells.Item(13, 4).Value = 4.701658921021388
$ws.Cells.Item(13, 5).Value = 7.383132016395942
$ws.Cells.Item(13, 6).Value = 32.18124386480401
$ws.Cells.Item(13, 7).Value = 42.02792245541253
$ws.Cells.Item(13, 9).Value = 3.655296337184813
$ws.Cells.Item(13, 10).Value = 13.32988563001441
$ws.Cells.Item(13, 11).Value = 21.87050815899161
$ws.Cells.Item(13, 12).Value = 5.866890006613863
$ws.Cells.Item(13, 13).Value = 14.19728239539849
$ws.Cells.Item(13, 14).Value = 7.534590252255419
$ws.Cells.Item(14, 2).Value = 17.02981682578852
$ws.Cells.Item(14, 3).Value = 8.121952473665903
$ws.Cells.Item(14, 4).Value = 4.675520689431637
$ws.Cells.Item(14, 5).Value = 7.369142996508341
$ws.Cells.Item(14, 6).Value = 32.08220710460298
$ws.Cells.Item(14, 7).Value = 41.87550737170691
$ws.Cells.Item(14, 9).Value = 3.70022243998498
$ws.Cells.Item(14, 10).Value = 13.31704358245627
$ws.Cells.Item(14, 11).Value = 21.84471597085402
$ws.Cells.Item(14, 12).Value = 5.870421845649449
$ws.Cells.Item(14, 13).Value = 14.07610737646825
$ws.Cells.Item(14, 14).Value = 7.508349954642123
$ws.Cells.Item(15, 2).Value = 16.93917298850139
$ws.Cells.Item(15, 3).Value = 8.088446026447803
$ws.Cells.Item(15, 4).Value = 4.661732170934266
$ws.Cells.Item(15, 5).Value = 7.36008056972701
$ws.Cells.Item(15, 6).Value = 32.00758422645498
$ws.Cells.Item(15, 7).Value = 41.7580609090102
$ws.Cells.Item(15, 9).Value = 3.730898398176099
$ws.Cells.Item(15, 10).Value = 13.30444729983996
$ws.Cells.Item(15, 11).Value = 21.81821103575616
$ws.Cells.Item(15, 12).Value = 5.871930348476522
$ws.Cells.Item(15, 13).Value = 14.00259775083024
$ws.Cells.Item(15, 14).Value = 7.494009151415558
$ws.Cells.Item(16, 2).Value = 16.42083199841958
$ws.Cells.Item(16, 3).Value = 7.860755204274053
$ws.Cells.Item(16, 4).Value = 4.568132293958498
$ws.Cells.Item(16, 5).Value = 7.311970668704337
$ws.Cells.Item(16, 6).Value = 31.67215566598178
$ws.Cells.Item(16, 7).Value = 41.24182973294081
$ws.Cells.Item(16, 9).Value = 3.888788667832936
$ws.Cells.Item(16, 10).Value = 13.26436451244334
$ws.Cells.Item(16, 11).Value = 21.73707760393223
$ws.Cells.Item(16, 12).Value = 5.885000018072605
$ws.Cells.Item(16, 13).Value = 13.56694530421974
$ws.Cells.Item(16, 14).Value = 7.401915075632558
$ws.Cells.Item(17, 2).Value = 16.09553272669619
$ws.Cells.Item(17, 3).Value = 7.719073795429618
$ws.Cells.Item(17, 4).Value = 4.510378722916419
$ws.Cells.Item(17, 5).Value = 7.282904350194702
$ws.Cells.Item(17, 6).Value = 31.46893458651346
$ws.Cells.Item(17, 7).Value = 40.92813610550357
$ws.Cells.Item(17, 9).Value = 3.98696461882672
$ws.Cells.Item(17, 10).Value = 13.24112171520205
$ws.Cells.Item(17, 11).Value = 21.68917692152183
$ws.Cells.Item(17, 12).Value = 5.893195670147207
$ws.Cells.Item(17, 13).Value = 13.29338676482499
$ws.Cells.Item(17, 14).Value = 7.346085399862454
$ws.Cells.Item(18, 2).Value = 15.9096706609274
$ws.Cells.Item(18, 3).Value = 7.622592003608177
$ws.Cells.Item(18, 4).Value = 4.471078716085079
$ws.Cells.Item(18, 5).Value = 7.267701178230797
$ws.Cells.Item(18, 6).Value = 31.38972227684699
$ws.Cells.Item(18, 7).Value = 40.81211539313509
$ws.Cells.Item(18, 9).Value = 4.035977612587646
$ws.Cells.Item(18, 10).Value = 13.24063966911076
$ws.Cells.Item(18, 11).Value = 21.69013602062414
$ws.Cells.Item(18, 12).Value = 5.899684676845132
$ws.Cells.Item(18, 13).Value = 13.13055363026432
$ws.Cells.Item(18, 14).Value = 7.309816273491869
$ws.Cells.Item(19, 2).Value = 15.84325748816939
$ws.Cells.Item(19, 3).Value = 7.600761678233272
$ws.Cells.Item(19, 4).Value = 4.462350171092905
$ws.Cells.Item(19, 5).Value = 7.261542445029878
$ws.Cells.Item(19, 6).Value = 31.33475653408919
$ws.Cells.Item(19, 7).Value = 40.72425117472071
$ws.Cells.Item(19, 9).Value = 4.058730932574258
$ws.Cells.Item(19, 10).Value = 13.23099570262262
$ws.Cells.Item(19, 11).Value = 21.66910867537092
$ws.Cells.Item(19, 12).Value = 5.900574458801276
$ws.Cells.Item(19, 13).Value = 13.07745500182786
$ws.Cells.Item(19, 14).Value = 7.300976867070706
$ws.Cells.Item(20, 2).Value = 16.1307445668718
$ws.Cells.Item(20, 3).Value = 7.733125994363349
$ws.Cells.Item(20, 4).Value = 4.516074169630633
$ws.Cells.Item(20, 5).Value = 7.286087043637446
$ws.Cells.Item(20, 6).Value = 31.49334026530487
$ws.Cells.Item(20, 7).Value = 40.96635670262842
$ws.Cells.Item(20, 9).Value = 3.975831045770145
$ws.Cells.Item(20, 10).Value = 13.24451677732774
$ws.Cells.Item(20, 11).Value = 21.69638579978937
$ws.Cells.Item(20, 12).Value = 5.89245037715947
$ws.Cells.Item(20, 13).Value = 13.3225027857908
$ws.Cells.Item(20, 14).Value = 7.351658548469553
$ws.Cells.Item(21, 2).Value = 17.06826650723109
$ws.Cells.Item(21, 3).Value = 8.155499809020318
$ws.Cells.Item(21, 4).Value = 4.68955276219873
$ws.Cells.Item(21, 5).Value = 7.371753607219716
$ws.Cells.Item(21, 6).Value = 32.07122168661646
$ws.Cells.Item(21, 7).Value = 41.85147427291317
$ws.Cells.Item(21, 9).Value = 3.695806305110584
$ws.Cells.Item(21, 10).Value = 13.30709399512463
$ws.Cells.Item(21, 11).Value = 21.82150641864622
$ws.Cells.Item(21, 12).Value = 5.867487381997071
$ws.Cells.Item(21, 13).Value = 14.1150904809536
$ws.Cells.Item(21, 14).Value = 7.520814654102646
$ws.Cells.Item(22, 2).Value = 17.69223915367544
$ws.Cells.Item(22, 3).Value = 8.406071542403339
$ws.Cells.Item(22, 4).Value = 4.793671292273948
$ws.Cells.Item(22, 5).Value = 7.431078292124314
$ws.Cells.Item(22, 6).Value = 32.50798081231869
$ws.Cells.Item(22, 7).Value = 42.52682379573636
$ws.Cells.Item(22, 9).Value = 3.505735862838888
$ws.Cells.Item(22, 10).Value = 13.36991026224558
$ws.Cells.Item(22, 11).Value = 21.94860916638606
$ws.Cells.Item(22, 12).Value = 5.854158030741906
$ws.Cells.Item(22, 13).Value = 14.60842902285207
$ws.Cells.Item(22, 14).Value = 7.627421391022566
$ws.Cells.Item(23, 2).Value = 17.3645910517651
$ws.Cells.Item(23, 3).Value = 8.255054491778251
$ws.Cells.Item(23, 4).Value = 4.730625250271374
$ws.Cells.Item(23, 5).Value = 7.401080481663594
$ws.Cells.Item(23, 6).Value = 32.32220886660738
$ws.Cells.Item(23, 7).Value = 42.24804331228723
$ws.Cells.Item(23, 9).Value = 3.596454095477925
$ws.Cells.Item(23, 10).Value = 13.35245323006417
$ws.Cells.Item(23, 11).Value = 21.91721436946884
$ws.Cells.Item(23, 12).Value = 5.863502018917811
$ws.Cells.Item(23, 13).Value = 14.34305753625497
$ws.Cells.Item(23, 14).Value = 7.564621182065723
$ws.Cells.Item(24, 2).Value = 16.12306635005104
$ws.Cells.Item(24, 3).Value = 7.696977968520752
$ws.Cells.Item(24, 4).Value = 4.500976619818594
$ws.Cells.Item(24, 5).Value = 7.287497022007755
$ws.Cells.Item(24, 6).Value = 31.559745485275
$ws.Cells.Item(24, 7).Value = 41.0824437293354
$ws.Cells.Item(24, 9).Value = 3.964665478526187
$ws.Cells.Item(24, 10).Value = 13.26919536155419
$ws.Cells.Item(24, 11).Value = 21.75198298081371
$ws.Cells.Item(24, 12).Value = 5.896403490311845
$ws.Cells.Item(24, 13).Value = 13.30273614142959
$ws.Cells.Item(24, 14).Value = 7.339860917651231
$ws.Cells.Item(25, 2).Value = 14.69546424886598
$ws.Cells.Item(25, 3).Value = 7.062753261653391
$ws.Cells.Item(25, 4).Value = 4.247744081447167
$ws.Cells.Item(25, 5).Value = 7.171510040726463
$ws.Cells.Item(25, 6).Value = 30.7791119943279
$ws.Cells.Item(25, 7).Value = 39.87820172120225
$ws.Cells.Item(25, 9).Value = 4.378953640874247
$ws.Cells.Item(25, 10).Value = 13.19929345030186
$ws.Cells.Item(25, 11).Value = 21.60264591935715
$ws.Cells.Item(25, 12).Value = 5.934604052869222
$ws.Cells.Item(25, 13).Value = 12.09061120884763
$ws.Cells.Item(25, 14).Value = 7.106353765727789
